$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.185.10"
$ws.Range("E2").Value = "  +1.99%  "
$ws.Range("D3").Value = "3.463.75"
$ws.Range("E3").Value = "  +1.38%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'580.90"
$ws.Range("E5").Value = "  +0.42%  "
$ws.Range("D6").Value = "'147.89"
$ws.Range("E6").Value = "  +2.33%  "
$ws.Range("D7").Value = "3.463.96"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +0.72%  "
$ws.Range("E10").Value = "  +2.47%  "
$ws.Range("E11").Value = "  +1.05%  "
$ws.Range("E12").Value = "  +5.20%  "
$ws.Range("D13").Value = "4.059.16"
$ws.Range("E13").Value = "  +1.47%  "
$ws.Range("D14").Value = "'29.36"
$ws.Range("E14").Value = "  +2.97%  "
$ws.Range("E15").Value = "  +2.42%  "
$ws.Range("D16").Value = "3.477.92"
$ws.Range("E16").Value = "  +1.62%  "
$ws.Range("E17").Value = "  +1.16%  "
$ws.Range("D18").Value = "63.208.53"
$ws.Range("E18").Value = "  +2.01%  "
$ws.Range("D19").Value = "'6.43"
$ws.Range("E19").Value = "  +4.16%  "
$ws.Range("D20").Value = "'14.50"
$ws.Range("E20").Value = "  +3.61%  "
$ws.Range("D21").Value = "'9.29"
$ws.Range("E21").Value = "  +1.62%  "
$ws.Range("D22").Value = "'388.44"
$ws.Range("E22").Value = "  -0.79%  "
$ws.Range("D23").Value = "'0.564"
$ws.Range("E23").Value = "  +1.95%  "
$ws.Range("D24").Value = "'74.52"
$ws.Range("E24").Value = "  -0.40%  "
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("D26").Value = "3.610.31"
$ws.Range("E26").Value = "  +1.49%  "
$ws.Range("E27").Value = "  +0.90%  "
$ws.Range("D28").Value = "'0.184"
$ws.Range("E28").Value = "  -0.23%  "
$ws.Range("D29").Value = "'7.64"
$ws.Range("E29").Value = "  +2.35%  "
$ws.Range("E30").Value = "  +0.09%  "
$ws.Range("D31").Value = "'8.18"
$ws.Range("E31").Value = "  +2.20%  "
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("D34").Value = "'23.44"
$ws.Range("E34").Value = "  -0.62%  "
$ws.Range("E35").Value = "  -4.29%  "
$ws.Range("D36").Value = "'5.33"
$ws.Range("E36").Value = "  +1.03%  "
$ws.Range("D37").Value = "'7.15"
$ws.Range("E37").Value = "  +2.50%  "
$ws.Range("D39").Value = "'31.87"
$ws.Range("E39").Value = "  +11.12%  "
$ws.Range("E40").Value = "  +0.46%  "
$ws.Range("D41").Value = "3.501.39"
$ws.Range("E41").Value = "  +1.66%  "
$ws.Range("D42").Value = "'0.0773"
$ws.Range("E42").Value = "  +2.33%  "
$ws.Range("D43").Value = "'0.793"
$ws.Range("E43").Value = "  +1.09%  "
$ws.Range("D44").Value = "'1.74"
$ws.Range("E44").Value = "  +3.94%  "
$ws.Range("E45").Value = "  -0.68%  "
$ws.Range("E46").Value = "  +3.80%  "
$ws.Range("D47").Value = "'4.37"
$ws.Range("E47").Value = "  -1.23%  "
$ws.Range("D48").Value = "2.592.45"
$ws.Range("E48").Value = "  +3.43%  "
$ws.Range("D49").Value = "'2.34"
$ws.Range("E49").Value = "  +11.69%  "
$ws.Range("D50").Value = "'6.82"
$ws.Range("E50").Value = "  +2.81%  "
$ws.Range("D51").Value = "'23.01"
$ws.Range("E51").Value = "  +0.72%  "
